$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values are plain text in this sheet (inline strings).
# Values that look like a simple decimal number (e.g. "612.14") get auto-coerced
# to a numeric type by the Value setter, so we force text via NumberFormat "@"
# first, then restore the default "Normal" style so no stray style is left behind.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.474.78"
$ws.Range("E2").Value = "  -0.35%  "
Set-TextValue $ws.Range("D3") "3.719.53"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue $ws.Range("D5") "612.14"
$ws.Range("E5").Value = "  -0.13%  "
Set-TextValue $ws.Range("D6") "178.15"
$ws.Range("E6").Value = "  +0.58%  "
Set-TextValue $ws.Range("D7") "3.718.67"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E9").Value = "  -2.66%  "
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("E11").Value = "  +2.23%  "
$ws.Range("E12").Value = "  -4.72%  "
Set-TextValue $ws.Range("D13") "39.49"
$ws.Range("E13").Value = "  -3.56%  "
Set-TextValue $ws.Range("D14") "0.0000252"
$ws.Range("E14").Value = "  -1.53%  "
Set-TextValue $ws.Range("D15") "4.335.74"
$ws.Range("E15").Value = "  -0.44%  "
Set-TextValue $ws.Range("D16") "3.714.03"
$ws.Range("E16").Value = "  -0.46%  "
Set-TextValue $ws.Range("D17") "69.515.42"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("E18").Value = "  -2.67%  "
Set-TextValue $ws.Range("D19") "7.46"
$ws.Range("E19").Value = "  -1.70%  "
Set-TextValue $ws.Range("D20") "499.95"
$ws.Range("E20").Value = "  -3.35%  "
Set-TextValue $ws.Range("D21") "16.25"
$ws.Range("E21").Value = "  -2.57%  "
$ws.Range("E22").Value = "  -2.30%  "
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("E24").Value = "  +3.46%  "
Set-TextValue $ws.Range("D25") "85.93"
$ws.Range("E25").Value = "  -2.85%  "
Set-TextValue $ws.Range("D26") "11.23"
$ws.Range("E26").Value = "  +2.72%  "
Set-TextValue $ws.Range("D27") "12.85"
$ws.Range("E27").Value = "  -5.33%  "
$ws.Range("E28").Value = "  +5.80%  "
$ws.Range("E30").Value = "  -2.95%  "
$ws.Range("E31").Value = "  +1.66%  "
Set-TextValue $ws.Range("D32") "7.99"
$ws.Range("E32").Value = "  +1.89%  "
Set-TextValue $ws.Range("D33") "30.26"
$ws.Range("E33").Value = "  -3.65%  "
$ws.Range("E34").Value = "  -2.42%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  +0.37%  "
Set-TextValue $ws.Range("D37") "6.05"
$ws.Range("E37").Value = "  -2.54%  "
Set-TextValue $ws.Range("D38") "0.346"
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("E39").Value = "  +3.79%  "
Set-TextValue $ws.Range("D40") "3.00"
$ws.Range("E40").Value = "  +10.22%  "
Set-TextValue $ws.Range("D41") "2.04"
$ws.Range("E41").Value = "  -6.80%  "
Set-TextValue $ws.Range("D42") "49.60"
$ws.Range("E42").Value = "  -3.29%  "
Set-TextValue $ws.Range("D43") "45.20"
$ws.Range("E43").Value = "  +1.03%  "
Set-TextValue $ws.Range("D44") "432.43"
$ws.Range("E44").Value = "  +2.13%  "
Set-TextValue $ws.Range("D45") "8.53"
$ws.Range("E45").Value = "  -3.59%  "
Set-TextValue $ws.Range("D46") "2.943.09"
$ws.Range("E46").Value = "  -3.73%  "
$ws.Range("E47").Value = "  -1.47%  "
Set-TextValue $ws.Range("D48") "139.42"
$ws.Range("E48").Value = "  +2.79%  "
Set-TextValue $ws.Range("D50") "26.88"
$ws.Range("E50").Value = "  -3.61%  "
Set-TextValue $ws.Range("D51") "2.45"
$ws.Range("E51").Value = "  -3.12%  "
